$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# "Add User Route and Model" - new Project Log Book entries for
# 09.03.2023 - 13.03.2023 (MERN / Node / Express / Mongo course work,
# User model + routes). 11 brand-new rows are inserted right after the
# last existing entry (row 332); together with the 10 already-blank
# filler rows below it, that gives 21 rows (333-353) of which 16 are
# filled in with data (333-348) and 5 stay blank (349-353), matching
# the target layout. The two "Total" rows (and everything after them)
# shift down by 11 automatically.
# ------------------------------------------------------------------

$ws.Rows("333:343").Insert() | Out-Null

# The newly inserted rows already inherit the A:E number formats /
# styles from the filler rows that used to sit there. Column F (the
# "Comment" column) is new for this block, so copy its formatting from
# an existing comment cell before filling in the text.
$ws.Range("F330").Copy() | Out-Null
$ws.Range("F332:F347").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 332
$ws.Range('A332').Value = '09.03.2023'
$ws.Range('B332').Value = 0.70833333333333337
$ws.Range('C332').Value = 'MERN Course - 1 '
$ws.Range('D332').Value = 'Tutorial'
$ws.Range('E332').Value = 45
$ws.Range('F332').Value = 'Take a Complete Node/Express/Mongo/TDD Course to be Able to Correctly Build the Backend of the App (Mosh Intro)'

# Row 333
$ws.Range('A333').Value = '09.03.2023'
$ws.Range('B333').Value = 0.76388888888888884
$ws.Range('C333').Value = 'MERN Course - 2'
$ws.Range('D333').Value = 'Tutorial'
$ws.Range('E333').Value = 25
$ws.Range('F333').Value = 'Node Architecture'

# Row 334
$ws.Range('A334').Value = '09.03.2023'
$ws.Range('B334').Value = 0.78125
$ws.Range('C334').Value = 'MERN Course - 3'
$ws.Range('D334').Value = 'Tutorial'
$ws.Range('E334').Value = 55
$ws.Range('F334').Value = 'Node Module System, Event Emitters, HTTP Module'

# Row 335
$ws.Range('A335').Value = '09.03.2023'
$ws.Range('B335').Value = 0.81944444444444453
$ws.Range('C335').Value = 'MERN Course - 4'
$ws.Range('D335').Value = 'Tutorial'
$ws.Range('E335').Value = 90
$ws.Range('F335').Value = 'NPM and RESTful Services with Express, Postman, HTTP Verbs, Joi'

# Row 336
$ws.Range('A336').Value = '10.03.2023'
$ws.Range('B336').Value = 0.84375
$ws.Range('C336').Value = 'MERN Course - 4 Cont'
$ws.Range('D336').Value = 'Tutorial'
$ws.Range('E336').Value = 20
$ws.Range('F336').Value = 'Continue Validation (Use Joi@13.1.0!!! Breaking Changes in Further Versions)'

# Row 337
$ws.Range('A337').Value = '10.03.2023'
$ws.Range('B337').Value = 0.85763888888888884
$ws.Range('C337').Value = 'Build User Request Sceleton'
$ws.Range('D337').Value = 'Code'
$ws.Range('E337').Value = 55
$ws.Range('F337').Value = 'Build the Sceleton Code for User GET, POST PUT and DELETE Requests'

# Row 338
$ws.Range('A338').Value = '10.03.2023'
$ws.Range('B338').Value = 0.92361111111111116
$ws.Range('C338').Value = 'MERN Course - 5'
$ws.Range('D338').Value = 'Tutorial'
$ws.Range('E338').Value = 55
$ws.Range('F338').Value = 'Advanced Express: Middlewares, Pipeline, Debug, ENV, Config'

# Row 339
$ws.Range('A339').Value = '10.03.2023'
$ws.Range('B339').Value = 0.96180555555555547
$ws.Range('C339').Value = 'Restructure App'
$ws.Range('D339').Value = 'Code'
$ws.Range('E339').Value = 60
$ws.Range('F339').Value = 'Apply Routes (Index, User)'

# Row 340
$ws.Range('A340').Value = '11.03.2023'
$ws.Range('B340').Value = 0.00694444444444444406
$ws.Range('C340').Value = 'MERN Course - 6'
$ws.Range('D340').Value = 'Tutorial'
$ws.Range('E340').Value = 75
$ws.Range('F340').Value = 'MongoDB, Mongoose'

# Row 341
$ws.Range('A341').Value = '11.03.2023'
$ws.Range('B341').Value = 0.05902777777777778317
$ws.Range('C341').Value = 'CRUD Users'
$ws.Range('D341').Value = 'Code'
$ws.Range('E341').Value = 60
$ws.Range('F341').Value = 'Create Functions for the Basic CRUD Operations for User Schema'

# Row 342
$ws.Range('A342').Value = '11.03.2023'
$ws.Range('B342').Value = 0.10069444444444443
$ws.Range('C342').Value = 'MERN Course - 7'
$ws.Range('D342').Value = 'Tutorial'
$ws.Range('E342').Value = 55
$ws.Range('F342').Value = 'MongoDB, Mongoose Data Validation'

# Row 343
$ws.Range('A343').Value = '11.03.2023'
$ws.Range('B343').Value = 0.1388888888888889
$ws.Range('C343').Value = 'User Model'
$ws.Range('D343').Value = 'Code'
$ws.Range('E343').Value = 25
$ws.Range('F343').Value = 'Create the User Model and Joi Validation'

# Row 344
$ws.Range('A344').Value = '11.03.2023'
$ws.Range('B344').Value = 0.16666666666666666
$ws.Range('C344').Value = 'Express Pipeline Paragraph'
$ws.Range('D344').Value = 'Documentation'
$ws.Range('E344').Value = 35
$ws.Range('F344').Value = 'Express Pipeline Explanation and Diagram'

# Row 345
$ws.Range('A345').Value = '11.03.2023'
$ws.Range('B345').Value = 0.19097222222222221
$ws.Range('C345').Value = 'Complete User Routes'
$ws.Range('D345').Value = 'Code'
$ws.Range('E345').Value = 45
$ws.Range('F345').Value = 'Complete and Test User Routes with MongoDB Connection (To be Continued...)'

# Row 346
$ws.Range('A346').Value = '11.03.2023'
$ws.Range('B346').Value = 0.55555555555555558
$ws.Range('C346').Value = 'Complete User Routes Continued'
$ws.Range('D346').Value = 'Code'
$ws.Range('E346').Value = 55
$ws.Range('F346').Value = 'Request GET/api/users/, GET/api/users/:id, POST/api/users{body}'

# Row 347
$ws.Range('A347').Value = '11.03.2023'
$ws.Range('B347').Value = 0.625
$ws.Range('C347').Value = 'Complete User Routes Continued'
$ws.Range('D347').Value = 'Code'
$ws.Range('E347').Value = 45
$ws.Range('F347').Value = 'Request PUT/api/users/id'

# Row 348
$ws.Range('A348').Value = '13.03.2023'
$ws.Range('B348').Value = 0.37152777777777773
$ws.Range('C348').Value = 'Profile Route'
$ws.Range('D348').Value = 'Code'

# Leave the selection where the author left it after typing the last
# new row.
$ws.Range("E348").Select() | Out-Null
